# Hazard interpretation template update:
# - Add a new "Negation words" column to the "Hazard-focused" sheet
#   (inserted as a new labeled column right after "Hazard words",
#   pushing the remaining columns two slots to the right, matching
#   the sheet's existing pattern of a labeled column followed by a
#   blank spacer column).
# - Make "Hazard-focused" the active/selected sheet tab, with the
#   new last header cell (M1) selected.

$wb = $excel.ActiveWorkbook

$wsHazard = $wb.Worksheets.Item("Hazard-focused")

# Insert two new columns at C:D on the Hazard-focused sheet. This shifts
# the existing headers (Hazard level 1 topics, Hazard level 2 topics,
# Best Documents, Hazard Category, Hazard name) from C/E/G/I/K to
# E/G/I/K/M respectively, preserving the blank spacer column pattern.
$wsHazard.Range("C1:D1").EntireColumn.Insert()

# Label the newly inserted column.
$wsHazard.Range("C1").Value = "Negation words"

# Switch the active tab to "Hazard-focused" and leave the cursor on the
# last header cell, mirroring where editing finished. (The "topic-focused"
# sheet's own selection, D12, is left untouched.)
$wsHazard.Activate()
$wsHazard.Range("M1").Select()
